$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Range("G2").Value = 6.713252999999999
    $ws.Range("H2").Value = 20.139759
    $ws.Range("I2").Value = 0.3101840064655811
    $ws.Range("J2").Value = 0.3231642354899327
    $ws.Range("M2").Value = 13.89934866666667
    $ws.Range("N2").Value = 41.69804600000001
    $ws.Range("O2").Value = 0.04853507553134179
    $ws.Range("P2").Value = 0.04999273878390351
    $ws.Range("Q2").Value = 93.309844134546
    $ws.Range("R2").Value = 839.788597210914
    $ws.Range("S2").Value = 0.01505480418242119
    $ws.Range("T2").Value = 0.01615586520914808
    $ws.Range("G3").Value = 6.713252999999999
    $ws.Range("H3").Value = 20.139759
    $ws.Range("I3").Value = 0.3101840064655811
    $ws.Range("J3").Value = 0.3231642354899327
    $ws.Range("O3").Value = 0.245697991654417
    $ws.Range("P3").Value = 0.253077086664408
    $ws.Range("Q3").Value = 472.3602684133049
    $ws.Range("R3").Value = 4251.242415719745
    $ws.Range("S3").Value = 0.07621158743191399
    $ws.Range("T3").Value = 0.08178546323192284
    $ws.Range("G4").Value = 6.713252999999999
    $ws.Range("H4").Value = 20.139759
    $ws.Range("I4").Value = 0.3101840064655811
    $ws.Range("J4").Value = 0.3231642354899327
    $ws.Range("M4").Value = 82.007665
    $ws.Range("N4").Value = 246.022995
    $ws.Range("O4").Value = 0.2863622109480123
    $ws.Range("P4").Value = 0.2949625822722868
    $ws.Range("Q4").Value = 550.538203084245
    $ws.Range("R4").Value = 4954.843827758204
    $ws.Range("S4").Value = 0.08882497789219637
    $ws.Range("T4").Value = 0.09532135739815993
    $ws.Range("G5").Value = 6.713252999999999
    $ws.Range("H5").Value = 20.139759
    $ws.Range("I5").Value = 0.3101840064655811
    $ws.Range("J5").Value = 0.3231642354899327
    $ws.Range("M5").Value = 25.0501465
    $ws.Range("N5").Value = 50.100293
    $ws.Range("O5").Value = 0.0874724982879541
    $ws.Range("P5").Value = 0.06006638442832619
    $ws.Range("Q5").Value = 168.1679711415645
    $ws.Range("R5").Value = 1009.007826849387
    $ws.Range("S5").Value = 0.02713256997451129
    $ws.Range("T5").Value = 0.01941130720242443
    $ws.Range("G6").Value = 6.713252999999999
    $ws.Range("H6").Value = 20.139759
    $ws.Range("I6").Value = 0.3101840064655811
    $ws.Range("J6").Value = 0.3231642354899327
    $ws.Range("M6").Value = 95.05788666666668
    $ws.Range("N6").Value = 285.17366
    $ws.Range("O6").Value = 0.3319322235782747
    $ws.Range("P6").Value = 0.3419012078510756
    $ws.Range("Q6").Value = 638.14764283866
    $ws.Range("R6").Value = 5743.32878554794
    $ws.Range("S6").Value = 0.1029600669845383
    $ws.Range("T6").Value = 0.1104902424482774
    $ws.Range("I7").Value = 0.0154484264788496
    $ws.Range("J7").Value = 0.01609489473505086
    $ws.Range("M7").Value = 13.89934866666667
    $ws.Range("N7").Value = 41.69804600000001
    $ws.Range("O7").Value = 0.04853507553134179
    $ws.Range("P7").Value = 0.04999273878390351
    $ws.Range("Q7").Value = 4.647210161770222
    $ws.Range("R7").Value = 41.824891455932
    $ws.Range("S7").Value = 0.0007497905459913456
    $ws.Range("T7").Value = 0.0008046278682438213
    $ws.Range("I8").Value = 0.0154484264788496
    $ws.Range("J8").Value = 0.01609489473505086
    $ws.Range("O8").Value = 0.245697991654417
    $ws.Range("P8").Value = 0.253077086664408
    $ws.Range("S8").Value = 0.003795647360074263
    $ws.Range("T8").Value = 0.004073249069716989
    $ws.Range("I9").Value = 0.0154484264788496
    $ws.Range("J9").Value = 0.01609489473505086
    $ws.Range("M9").Value = 82.007665
    $ws.Range("N9").Value = 246.022995
    $ws.Range("O9").Value = 0.2863622109480123
    $ws.Range("P9").Value = 0.2949625822722868
    $ws.Range("Q9").Value = 27.41904410564333
    $ws.Range("R9").Value = 246.77139695079
    $ws.Range("S9").Value = 0.004423845562151188
    $ws.Range("T9").Value = 0.004747391712451234
    $ws.Range("I10").Value = 0.0154484264788496
    $ws.Range("J10").Value = 0.01609489473505086
    $ws.Range("M10").Value = 25.0501465
    $ws.Range("N10").Value = 50.100293
    $ws.Range("O10").Value = 0.0874724982879541
    $ws.Range("P10").Value = 0.06006638442832619
    $ws.Range("Q10").Value = 8.375449681884334
    $ws.Range("R10").Value = 50.252698091306
    $ws.Range("S10").Value = 0.001351312458722756
    $ws.Range("T10").Value = 0.0009667621344890078
    $ws.Range("I11").Value = 0.0154484264788496
    $ws.Range("J11").Value = 0.01609489473505086
    $ws.Range("M11").Value = 95.05788666666668
    $ws.Range("N11").Value = 285.17366
    $ws.Range("O11").Value = 0.3319322235782747
    $ws.Range("P11").Value = 0.3419012078510756
    $ws.Range("Q11").Value = 31.78235091930222
    $ws.Range("R11").Value = 286.04115827372
    $ws.Range("S11").Value = 0.005127830551910043
    $ws.Range("T11").Value = 0.005502863950149805
    $ws.Range("G12").Value = 6.661784666666667
    $ws.Range("H12").Value = 19.985354
    $ws.Range("I12").Value = 0.3078059262949933
    $ws.Range("J12").Value = 0.3206866401135023
    $ws.Range("M12").Value = 13.89934866666667
    $ws.Range("N12").Value = 41.69804600000001
    $ws.Range("O12").Value = 0.04853507553134179
    $ws.Range("P12").Value = 0.04999273878390351
    $ws.Range("Q12").Value = 92.59446782425378
    $ws.Range("R12").Value = 833.3502104182842
    $ws.Range("S12").Value = 0.01493938388172212
    $ws.Range("T12").Value = 0.01603200343068199
    $ws.Range("G13").Value = 6.661784666666667
    $ws.Range("H13").Value = 19.985354
    $ws.Range("I13").Value = 0.3078059262949933
    $ws.Range("J13").Value = 0.3206866401135023
    $ws.Range("O13").Value = 0.245697991654417
    $ws.Range("P13").Value = 0.253077086664408
    $ws.Range("Q13").Value = 468.7388354436078
    $ws.Range("R13").Value = 4218.64951899247
    $ws.Range("S13").Value = 0.07562729791000737
    $ws.Range("T13").Value = 0.08115844061212263
    $ws.Range("G14").Value = 6.661784666666667
    $ws.Range("H14").Value = 19.985354
    $ws.Range("I14").Value = 0.3078059262949933
    $ws.Range("J14").Value = 0.3206866401135023
    $ws.Range("M14").Value = 82.007665
    $ws.Range("N14").Value = 246.022995
    $ws.Range("O14").Value = 0.2863622109480123
    $ws.Range("P14").Value = 0.2949625822722868
    $ws.Range("Q14").Value = 546.3174052461367
    $ws.Range("R14").Value = 4916.856647215231
    $ws.Range("S14").Value = 0.08814398559673521
    $ws.Range("T14").Value = 0.09459055946810214
    $ws.Range("G15").Value = 6.661784666666667
    $ws.Range("H15").Value = 19.985354
    $ws.Range("I15").Value = 0.3078059262949933
    $ws.Range("J15").Value = 0.3206866401135023
    $ws.Range("M15").Value = 25.0501465
    $ws.Range("N15").Value = 50.100293
    $ws.Range("O15").Value = 0.0874724982879541
    $ws.Range("P15").Value = 0.06006638442832619
    $ws.Range("Q15").Value = 166.8786818514537
    $ws.Range("R15").Value = 1001.272091108722
    $ws.Range("S15").Value = 0.02692455336086093
    $ws.Range("T15").Value = 0.01926248700608592
    $ws.Range("G16").Value = 6.661784666666667
    $ws.Range("H16").Value = 19.985354
    $ws.Range("I16").Value = 0.3078059262949933
    $ws.Range("J16").Value = 0.3206866401135023
    $ws.Range("M16").Value = 95.05788666666668
    $ws.Range("N16").Value = 285.17366
    $ws.Range("O16").Value = 0.3319322235782747
    $ws.Range("P16").Value = 0.3419012078510756
    $ws.Range("Q16").Value = 633.2551718417378
    $ws.Range("R16").Value = 5699.296546575642
    $ws.Range("S16").Value = 0.1021707055456677
    $ws.Range("T16").Value = 0.1096431495965096
    $ws.Range("G17").Value = 2.607918
    $ws.Range("H17").Value = 5.215835999999999
    $ws.Range("I17").Value = 0.1204981331366039
    $ws.Range("J17").Value = 0.08369373503331734
    $ws.Range("M17").Value = 13.89934866666667
    $ws.Range("N17").Value = 41.69804600000001
    $ws.Range("O17").Value = 0.04853507553134179
    $ws.Range("P17").Value = 0.04999273878390351
    $ws.Range("Q17").Value = 36.248361576076
    $ws.Range("R17").Value = 217.490169456456
    $ws.Range("S17").Value = 0.005848385993170748
    $ws.Range("T17").Value = 0.004184079033369867
    $ws.Range("G18").Value = 2.607918
    $ws.Range("H18").Value = 5.215835999999999
    $ws.Range("I18").Value = 0.1204981331366039
    $ws.Range("J18").Value = 0.08369373503331734
    $ws.Range("O18").Value = 0.245697991654417
    $ws.Range("P18").Value = 0.253077086664408
    $ws.Range("Q18").Value = 183.49924343383
    $ws.Range("R18").Value = 1100.99546060298
    $ws.Range("S18").Value = 0.02960614930977013
    $ws.Range("T18").Value = 0.02118096663429485
    $ws.Range("G19").Value = 2.607918
    $ws.Range("H19").Value = 5.215835999999999
    $ws.Range("I19").Value = 0.1204981331366039
    $ws.Range("J19").Value = 0.08369373503331734
    $ws.Range("M19").Value = 82.007665
    $ws.Range("N19").Value = 246.022995
    $ws.Range("O19").Value = 0.2863622109480123
    $ws.Range("P19").Value = 0.2949625822722868
    $ws.Range("Q19").Value = 213.86926569147
    $ws.Range("R19").Value = 1283.21559414882
    $ws.Range("S19").Value = 0.03450611182010584
    $ws.Range("T19").Value = 0.02468652020543984
    $ws.Range("G20").Value = 2.607918
    $ws.Range("H20").Value = 5.215835999999999
    $ws.Range("I20").Value = 0.1204981331366039
    $ws.Range("J20").Value = 0.08369373503331734
    $ws.Range("M20").Value = 25.0501465
    $ws.Range("N20").Value = 50.100293
    $ws.Range("O20").Value = 0.0874724982879541
    $ws.Range("P20").Value = 0.06006638442832619
    $ws.Range("Q20").Value = 65.32872795998699
    $ws.Range("R20").Value = 261.314911839948
    $ws.Range("S20").Value = 0.01054027274449325
    $ws.Range("T20").Value = 0.00502718006275371
    $ws.Range("G21").Value = 2.607918
    $ws.Range("H21").Value = 5.215835999999999
    $ws.Range("I21").Value = 0.1204981331366039
    $ws.Range("J21").Value = 0.08369373503331734
    $ws.Range("M21").Value = 95.05788666666668
    $ws.Range("N21").Value = 285.17366
    $ws.Range("O21").Value = 0.3319322235782747
    $ws.Range("P21").Value = 0.3419012078510756
    $ws.Range("Q21").Value = 247.90317367996
    $ws.Range("R21").Value = 1487.41904207976
    $ws.Range("S21").Value = 0.03999721326906391
    $ws.Range("T21").Value = 0.02861498909745908
    $ws.Range("G22").Value = 5.325505333333333
    $ws.Range("H22").Value = 15.976516
    $ws.Range("I22").Value = 0.2460635076239721
    $ws.Range("J22").Value = 0.2563604946281968
    $ws.Range("M22").Value = 13.89934866666667
    $ws.Range("N22").Value = 41.69804600000001
    $ws.Range("O22").Value = 0.04853507553134179
    $ws.Range("P22").Value = 0.04999273878390351
    $ws.Range("Q22").Value = 74.02105545419289
    $ws.Range("R22").Value = 666.1894990877361
    $ws.Range("S22").Value = 0.01194271092803638
    $ws.Range("T22").Value = 0.01281616324245974
    $ws.Range("G23").Value = 5.325505333333333
    $ws.Range("H23").Value = 15.976516
    $ws.Range("I23").Value = 0.2460635076239721
    $ws.Range("J23").Value = 0.2563604946281968
    $ws.Range("O23").Value = 0.245697991654417
    $ws.Range("P23").Value = 0.253077086664408
    $ws.Range("Q23").Value = 374.7150790667089
    $ws.Range("R23").Value = 3372.43571160038
    $ws.Range("S23").Value = 0.06045730964265128
    $ws.Range("T23").Value = 0.06487896711635065
    $ws.Range("G24").Value = 5.325505333333333
    $ws.Range("H24").Value = 15.976516
    $ws.Range("I24").Value = 0.2460635076239721
    $ws.Range("J24").Value = 0.2563604946281968
    $ws.Range("M24").Value = 82.007665
    $ws.Range("N24").Value = 246.022995
    $ws.Range("O24").Value = 0.2863622109480123
    $ws.Range("P24").Value = 0.2949625822722868
    $ws.Range("Q24").Value = 436.7322573317134
    $ws.Range("R24").Value = 3930.59031598542
    $ws.Range("S24").Value = 0.07046329007682374
    $ws.Range("T24").Value = 0.07561675348813363
    $ws.Range("G25").Value = 5.325505333333333
    $ws.Range("H25").Value = 15.976516
    $ws.Range("I25").Value = 0.2460635076239721
    $ws.Range("J25").Value = 0.2563604946281968
    $ws.Range("M25").Value = 25.0501465
    $ws.Range("N25").Value = 50.100293
    $ws.Range("O25").Value = 0.0874724982879541
    $ws.Range("P25").Value = 0.06006638442832619
    $ws.Range("Q25").Value = 133.4046887865313
    $ws.Range("R25").Value = 800.4281327191881
    $ws.Range("S25").Value = 0.02152378974936588
    $ws.Range("T25").Value = 0.01539864802257312
    $ws.Range("G26").Value = 5.325505333333333
    $ws.Range("H26").Value = 15.976516
    $ws.Range("I26").Value = 0.2460635076239721
    $ws.Range("J26").Value = 0.2563604946281968
    $ws.Range("M26").Value = 95.05788666666668
    $ws.Range("N26").Value = 285.17366
    $ws.Range("O26").Value = 0.3319322235782747
    $ws.Range("P26").Value = 0.3419012078510756
    $ws.Range("Q26").Value = 506.2312824187289
    $ws.Range("R26").Value = 4556.081541768561
    $ws.Range("S26").Value = 0.08167640722709481
    $ws.Range("T26").Value = 0.08764996275867966
